# Daily cryptos-list refresh (GitHub Actions scheduled update).
# Updates Price (col D) and Volume(1h) (col E) text values for existing
# coins, and replaces row 51 (Bittensor -> InjectiveProtocol) entirely.
# Leading "'" forces text entry for values that would otherwise be
# auto-parsed as numbers by Excel, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.593.33'
$ws.Range("E2").Value = '  -1.63%  '

$ws.Range("D3").Value = '2.898.14'
$ws.Range("E3").Value = '  -2.94%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''529.19'
$ws.Range("E5").Value = '  -2.80%  '

$ws.Range("D6").Value = '''143.39'
$ws.Range("E6").Value = '  -6.29%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").Value = '''0.557'
$ws.Range("E8").Value = '  -1.07%  '

$ws.Range("D9").Value = '2.903.47'
$ws.Range("E9").Value = '  -2.75%  '

$ws.Range("E10").Value = '  -2.99%  '

$ws.Range("D11").Value = '''5.96'
$ws.Range("E11").Value = '  -3.87%  '

$ws.Range("D12").Value = '''0.363'
$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("D13").Value = '3.411.30'
$ws.Range("E13").Value = '  -2.60%  '

$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").Value = '60.513.49'
$ws.Range("E15").Value = '  -1.90%  '

$ws.Range("D16").Value = '''22.70'
$ws.Range("E16").Value = '  -3.96%  '

$ws.Range("D17").Value = '2.899.56'
$ws.Range("E17").Value = '  -2.78%  '

$ws.Range("D18").Value = '''0.0000143'
$ws.Range("E18").Value = '  -3.07%  '

$ws.Range("D19").Value = '''5.05'
$ws.Range("E19").Value = '  -1.24%  '

$ws.Range("D20").Value = '''11.74'
$ws.Range("E20").Value = '  -1.31%  '

$ws.Range("D21").Value = '''364.72'
$ws.Range("E21").Value = '  -6.10%  '

$ws.Range("D22").Value = '''6.66'
$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '''64.26'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").Value = '3.018.99'
$ws.Range("E25").Value = '  -2.84%  '

$ws.Range("D26").Value = '''0.453'
$ws.Range("E26").Value = '  -3.32%  '

$ws.Range("D27").Value = '''0.179'
$ws.Range("E27").Value = '  -4.59%  '

$ws.Range("D28").Value = '''0.999'

$ws.Range("D29").Value = '''7.81'
$ws.Range("E29").Value = '  -6.61%  '

$ws.Range("D30").Value = '0.0₃0865'
$ws.Range("E30").Value = '  -8.03%  '

$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("E32").Value = '  -2.56%  '

$ws.Range("D33").Value = '''19.71'
$ws.Range("E33").Value = '  -2.82%  '

$ws.Range("D34").Value = '''144.38'
$ws.Range("E34").Value = '  -9.11%  '

$ws.Range("D35").Value = '''4.36'
$ws.Range("E35").Value = '  -5.29%  '

$ws.Range("D36").Value = '''5.59'
$ws.Range("E36").Value = '  -7.04%  '

$ws.Range("D37").Value = '''0.998'
$ws.Range("E37").Value = '  -6.63%  '

$ws.Range("D38").Value = '''1.21'
$ws.Range("E38").Value = '  -5.64%  '

$ws.Range("D39").Value = '''37.64'
$ws.Range("E39").Value = '  +1.25%  '

$ws.Range("D40").Value = '''1.50'
$ws.Range("E40").Value = '  -5.31%  '

$ws.Range("D41").Value = '2.323.77'
$ws.Range("E41").Value = '  -4.81%  '

$ws.Range("D42").Value = '''3.69'
$ws.Range("E42").Value = '  -4.74%  '

$ws.Range("D43").Value = '''0.645'
$ws.Range("E43").Value = '  -2.19%  '

$ws.Range("D44").Value = '''0.0581'
$ws.Range("E44").Value = '  -2.25%  '

$ws.Range("D45").Value = '''20.71'
$ws.Range("E45").Value = '  -6.82%  '

$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = '''4.93'
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("D48").Value = '''0.0235'
$ws.Range("E48").Value = '  -4.28%  '

$ws.Range("D49").Value = '''0.0936'
$ws.Range("E49").Value = '  -2.25%  '

$ws.Range("D50").Value = '''10.32'
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''18.53'
$ws.Range("E51").Value = '  -6.01%  '
